$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    "B2" = 1.02
    "C2" = 1.071078228498918
    "D2" = 1.061183748586212
    "E2" = 1.084771117061936
    "F2" = 1.091975180185917
    "I2" = 1.041242387694905
    "J2" = 1.076004643915476
    "K2" = 1.063908500353023
    "L2" = 1.087433274228871
    "M2" = 1.094618801330121
    "N2" = 1.077532693929294
    "B3" = 1.02
    "C3" = 1.073225611200939
    "D3" = 1.062748287427633
    "E3" = 1.086955617671272
    "F3" = 1.094309064427237
    "I3" = 1.041716542126218
    "J3" = 1.07780420592111
    "K3" = 1.065285720743031
    "L3" = 1.08943356105254
    "M3" = 1.096769496314435
    "N3" = 1.079334811519388
    "B4" = 1.02
    "C4" = 1.07460943965118
    "D4" = 1.063755611417724
    "E4" = 1.088363775259226
    "F4" = 1.095813968642354
    "I4" = 1.042019819090874
    "J4" = 1.078962812116172
    "K4" = 1.066171353673845
    "L4" = 1.090722151676594
    "M4" = 1.098155521383163
    "N4" = 1.080495063068142
    "B5" = 1.02
    "C5" = 1.075189875220028
    "D5" = 1.064177906796334
    "E5" = 1.08895451302169
    "F5" = 1.096445401165437
    "I5" = 1.04214647933247
    "J5" = 1.079448521541835
    "K5" = 1.066542373259453
    "L5" = 1.091262533476441
    "M5" = 1.098736892914698
    "N5" = 1.080981462256901
    "B6" = 1.02
    "C6" = 1.075287256076027
    "D6" = 1.064248743284328
    "E6" = 1.089053628004692
    "F6" = 1.096551350354559
    "I6" = 1.042167697279688
    "J6" = 1.079529994823058
    "K6" = 1.066604593430315
    "L6" = 1.091353188140443
    "M6" = 1.09883443164847
    "N6" = 1.081063051239528
    "B7" = 1.02
    "C7" = 1.074617200626703
    "D7" = 1.063761258770047
    "E7" = 1.088371673601929
    "F7" = 1.095822410659405
    "I7" = 1.04202151481036
    "J7" = 1.07896930753035
    "K7" = 1.066176316337261
    "L7" = 1.090729377516599
    "M7" = 1.098163294816284
    "N7" = 1.080501567706552
    "B8" = 1.02
    "C8" = 1.071805140248714
    "D8" = 1.061713548327076
    "E8" = 1.085510508318006
    "F8" = 1.092765039878578
    "I8" = 1.041403366916429
    "J8" = 1.076614038145928
    "K8" = 1.064375094377888
    "L8" = 1.088110484728601
    "M8" = 1.095346821272433
    "N8" = 1.078142953569452
    "B9" = 1.02
    "C9" = 1.066805045726392
    "D9" = 1.058065630383912
    "E9" = 1.080426283848626
    "F9" = 1.087335638874361
    "I9" = 1.040286677261447
    "J9" = 1.072417859781438
    "K9" = 1.061157868392248
    "L9" = 1.083450437239441
    "M9" = 1.090339370901535
    "N9" = 1.073940816150376
    "B10" = 1.02
    "C10" = 1.06343943207196
    "D10" = 1.055605639302639
    "E10" = 1.07700621698209
    "F10" = 1.083685703655597
    "I10" = 1.039523237706534
    "J10" = 1.069587807784158
    "K10" = 1.058982621264047
    "L10" = 1.080311430537017
    "M10" = 1.086969153212357
    "N10" = 1.071106745154648
    "B11" = 1.02
    "C11" = 1.061973983880016
    "D11" = 1.054533462367176
    "E11" = 1.075517581944003
    "F11" = 1.082097563310592
    "I11" = 1.039188037489322
    "J11" = 1.068354239195502
    "K11" = 1.05803319223013
    "L11" = 1.078944114725592
    "M11" = 1.085501786346033
    "N11" = 1.069871424757104
    "B12" = 1.02
    "C12" = 1.061428392108094
    "D12" = 1.054134130795708
    "E12" = 1.07496343679943
    "F12" = 1.081506459355087
    "I12" = 1.03906282368299
    "J12" = 1.067894780227683
    "K12" = 1.057679373736917
    "L12" = 1.078434977851239
    "M12" = 1.084955493473697
    "N12" = 1.06941131330485
    "B13" = 1.02
    "C13" = 1.061545481017363
    "D13" = 1.054219837978475
    "E13" = 1.07508235774832
    "F13" = 1.081633308049371
    "I13" = 1.039089714566998
    "J13" = 1.067993393156873
    "K13" = 1.057755321850849
    "L13" = 1.078544246829898
    "M13" = 1.085072732242021
    "N13" = 1.069510066275711
    "B14" = 1.02
    "C14" = 1.061928911011441
    "D14" = 1.05450047563546
    "E14" = 1.075471800839876
    "F14" = 1.082048727155657
    "I14" = 1.0391777017311
    "J14" = 1.068316286009095
    "K14" = 1.058003969267053
    "L14" = 1.078902055129835
    "M14" = 1.085456655286123
    "N14" = 1.069833417672819
    "B15" = 1.02
    "C15" = 1.062164986843628
    "D15" = 1.054673242157746
    "E15" = 1.075711589403691
    "F15" = 1.082304520574541
    "I15" = 1.039231819719906
    "J15" = 1.06851506340032
    "K15" = 1.058157014860463
    "L15" = 1.079122345222155
    "M15" = 1.085693036447145
    "N15" = 1.070032477350739
    "B16" = 1.02
    "C16" = 1.063536514678114
    "D16" = 1.055676646622007
    "E16" = 1.077104846710448
    "F16" = 1.083790937626597
    "I16" = 1.039545385484385
    "J16" = 1.069669501324525
    "K16" = 1.059045470779491
    "L16" = 1.080402000857588
    "M16" = 1.087066364760546
    "N16" = 1.071188554709213
    "B17" = 1.02
    "C17" = 1.064394637258281
    "D17" = 1.056304166177033
    "E17" = 1.077976704776622
    "F17" = 1.08472123840924
    "I17" = 1.039740831151798
    "J17" = 1.070391447222897
    "K17" = 1.059600740378285
    "L17" = 1.081202498987691
    "M17" = 1.08792563701658
    "N17" = 1.071911525853568
    "B18" = 1.02
    "C18" = 1.064894385359911
    "D18" = 1.056669515721171
    "E18" = 1.078484502326232
    "F18" = 1.085263127595278
    "I18" = 1.039854385696796
    "J18" = 1.070811763244026
    "K18" = 1.05992389476571
    "L18" = 1.081668636513631
    "M18" = 1.088426062909397
    "N18" = 1.072332438771659
    "B19" = 1.02
    "C19" = 1.065064655394115
    "D19" = 1.05679397733116
    "E19" = 1.078657523408926
    "F19" = 1.085447773626905
    "I19" = 1.039893029595638
    "J19" = 1.070954948391306
    "K19" = 1.060033959889719
    "L19" = 1.081827446264755
    "M19" = 1.088596565265087
    "N19" = 1.072475827258272
    "B20" = 1.02
    "C20" = 1.064302649821079
    "D20" = 1.056236908976246
    "E20" = 1.077883239742042
    "F20" = 1.084621502733092
    "I20" = 1.039719907848246
    "J20" = 1.070314070475051
    "K20" = 1.059541240383963
    "L20" = 1.081116694007769
    "M20" = 1.087833525482473
    "N20" = 1.071834039221865
    "B21" = 1.02
    "C21" = 1.061816035514609
    "D21" = 1.054417864806501
    "E21" = 1.075357152996768
    "F21" = 1.081926429991723
    "I21" = 1.039151811265703
    "J21" = 1.068221237079171
    "K21" = 1.057930780990223
    "L21" = 1.078796724446941
    "M21" = 1.085343634276009
    "N21" = 1.069738233762512
    "B22" = 1.02
    "C22" = 1.060245297610713
    "D22" = 1.053267912291113
    "E22" = 1.073761940234146
    "F22" = 1.080224978696507
    "I22" = 1.038790539110787
    "J22" = 1.066898101622582
    "K22" = 1.056911506310972
    "L22" = 1.077330788133676
    "M22" = 1.083770903223467
    "N22" = 1.068413219301808
    "B23" = 1.02
    "C23" = 1.06107868138452
    "D23" = 1.053878125510489
    "E23" = 1.07460826594355
    "F23" = 1.081127623390382
    "I23" = 1.038982447311518
    "J23" = 1.067600223130985
    "K23" = 1.057452488598669
    "L23" = 1.078108611295494
    "M23" = 1.084605337125529
    "N23" = 1.069116337903288
    "B24" = 1.02
    "C24" = 1.06434421739387
    "D24" = 1.056267301685646
    "E24" = 1.077925474860812
    "F24" = 1.084666571275527
    "I24" = 1.039729363565009
    "J24" = 1.070349036114279
    "K24" = 1.05956812810939
    "L24" = 1.081155467991911
    "M24" = 1.087875149107076
    "N24" = 1.071869054516311
    "B25" = 1.02
    "C25" = 1.068103216999295
    "D25" = 1.059013541773507
    "E25" = 1.081745918342872
    "F25" = 1.088744456257015
    "I25" = 1.040578673853584
    "J25" = 1.07350828878378
    "K25" = 1.061994861098169
    "L25" = 1.084660724411812
    "M25" = 1.091639392246854
    "N25" = 1.075032793686974
}

foreach ($addr in $values.Keys) {
    $ws.Range($addr).Value = $values[$addr]
}
